$d = $word.ActiveDocument
$d.Content.Find.Execute("CS-XXXX.YYY", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Dr. Bastani", 2)
